$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.241.64"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "4.004.91"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'530.33"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'149.39"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'0.736"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "'0.175"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'0.0000343"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").Value = "'44.31"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "'10.64"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "4.647.40"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "4.019.80"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "'21.48"
$ws.Range("E16").Value = "  +7.45%  "
$ws.Range("D17").Value = "'14.32"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'0.134"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "71.229.02"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").Value = "'440.50"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'3.56"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "'93.57"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'14.40"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'12.41"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("D26").Value = "'4.12"
$ws.Range("E26").Value = "  +5.62%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").Value = "'36.91"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "'699.39"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D33").Value = "'6.90"
$ws.Range("E33").Value = "  +13.91%  "
$ws.Range("D34").Value = "'66.67"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "0.0₃0907"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").Value = "'0.442"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").Value = "'41.17"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'3.54"
$ws.Range("E38").Value = "  +15.62%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "'2.94"
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'3.52"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "  +7.93%  "
$ws.Range("D47").Value = "'0.146"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'0.000285"
$ws.Range("E48").Value = "  +20.13%  "
$ws.Range("D49").Value = "'9.29"
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -3.58%  "
